# chore: update Sheets via scheduled runner
# Refresh of market-price derived columns (currentAveragePrice[NQ/HQ],
# LevePrice[NQ/HQ], LeveProfit[NQ/HQ]) on a handful of rows across the
# ALC / ARM / BSM / CUL / GSM / LTW / WVR leve-profit sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 766.73334
$ws.Range("I8").Value = 59.18182
$ws.Range("J8").Value = 2712.5
$ws.Range("K8").Value = 177.54546
$ws.Range("L8").Value = 8137.5
$ws.Range("M8").Value = -38.54545999999999
$ws.Range("N8").Value = -8415.5

$ws.Range("H64").Value = 3072.375
$ws.Range("I64").Value = 2986.0625
$ws.Range("J64").Value = 3245
$ws.Range("K64").Value = 2986.0625
$ws.Range("L64").Value = 3245
$ws.Range("M64").Value = -2738.0625
$ws.Range("N64").Value = -3741

$ws.Range("H67").Value = 3072.375
$ws.Range("I67").Value = 2986.0625
$ws.Range("J67").Value = 3245
$ws.Range("K67").Value = 2986.0625
$ws.Range("L67").Value = 3245
$ws.Range("M67").Value = -2128.0625
$ws.Range("N67").Value = -4961

$ws.Range("H116").Value = 2916.4473
$ws.Range("I116").Value = 2719.625
$ws.Range("J116").Value = 3059.5908
$ws.Range("K116").Value = 2719.625
$ws.Range("L116").Value = 3059.5908
$ws.Range("M116").Value = 722.375
$ws.Range("N116").Value = -9943.5908

$ws.Range("H132").Value = 3637723
$ws.Range("I132").Value = 4082522.5
$ws.Range("K132").Value = 12247567.5
$ws.Range("M132").Value = -12245037.5

$ws.Range("H135").Value = 798.1053000000001
$ws.Range("I135").Value = 829.3333
$ws.Range("K135").Value = 7463.9997
$ws.Range("M135").Value = -4928.9997

$ws.Range("H137").Value = 2227.8696
$ws.Range("I137").Value = 2158.3235
$ws.Range("J137").Value = 2424.9167
$ws.Range("K137").Value = 6474.970499999999
$ws.Range("L137").Value = 7274.750100000001
$ws.Range("M137").Value = -3924.970499999999
$ws.Range("N137").Value = -12374.7501

$ws.Range("H138").Value = 4431.6514
$ws.Range("I138").Value = 2016.5483
$ws.Range("J138").Value = 10670.667
$ws.Range("K138").Value = 6049.644899999999
$ws.Range("L138").Value = 32012.001
$ws.Range("M138").Value = -909.6448999999993
$ws.Range("N138").Value = -42292.001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").ClearContents()
$ws.Range("N24").Value = 0

$ws.Range("H32").Value = 5069.772
$ws.Range("I32").Value = 3783.7073
$ws.Range("J32").Value = 8365.3125
$ws.Range("K32").Value = 3783.7073
$ws.Range("L32").Value = 8365.3125
$ws.Range("M32").Value = -3496.7073
$ws.Range("N32").Value = -8939.3125

$ws.Range("H33").Value = 17718

$ws.Range("H36").Value = 25939.572
$ws.Range("I36").Value = 8303.799999999999
$ws.Range("K36").Value = 8303.799999999999
$ws.Range("M36").Value = -7957.799999999999

$ws.Range("H61").Value = 3284.0967
$ws.Range("I61").Value = 1725.3636
$ws.Range("J61").Value = 4141.4
$ws.Range("K61").Value = 1725.3636
$ws.Range("L61").Value = 4141.4
$ws.Range("M61").Value = -1513.3636
$ws.Range("N61").Value = -4565.4

$ws.Range("H94").Value = 30000
$ws.Range("J94").Value = 30000
$ws.Range("L94").Value = 30000
$ws.Range("N94").Value = -31802

$ws.Range("H97").Value = 581.4706
$ws.Range("I97").Value = 633.53845
$ws.Range("J97").Value = 412.25
$ws.Range("K97").Value = 633.53845
$ws.Range("L97").Value = 412.25
$ws.Range("M97").Value = -137.53845
$ws.Range("N97").Value = -1404.25

$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").ClearContents()
$ws.Range("N100").Value = 0

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").ClearContents()
$ws.Range("N109").Value = 0

$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").ClearContents()
$ws.Range("N112").Value = 0

$ws.Range("H115").Value = 29975
$ws.Range("J115").Value = 29975
$ws.Range("L115").Value = 29975
$ws.Range("N115").Value = -33109

$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").ClearContents()
$ws.Range("N121").Value = 0

$ws.Range("H136").Value = 3284.0967
$ws.Range("I136").Value = 1725.3636
$ws.Range("J136").Value = 4141.4
$ws.Range("K136").Value = 5176.0908
$ws.Range("L136").Value = 12424.2
$ws.Range("M136").Value = -2626.0908
$ws.Range("N136").Value = -17524.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1708.8572
$ws.Range("I105").Value = 1444.5454
$ws.Range("J105").Value = 2678
$ws.Range("K105").Value = 1444.5454
$ws.Range("L105").Value = 2678
$ws.Range("M105").Value = 302.4546
$ws.Range("N105").Value = -6172

$ws.Range("H134").Value = 2991.3635
$ws.Range("I134").Value = 2072.7778
$ws.Range("J134").Value = 7125
$ws.Range("K134").Value = 6218.3334
$ws.Range("L134").Value = 21375
$ws.Range("M134").Value = -3683.3334
$ws.Range("N134").Value = -26445

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 3117.6667
$ws.Range("I132").Value = 1676.75
$ws.Range("J132").Value = 3838.125
$ws.Range("K132").Value = 15090.75
$ws.Range("L132").Value = 34543.125
$ws.Range("M132").Value = -12560.75
$ws.Range("N132").Value = -39603.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2330.1177
$ws.Range("I102").Value = 1589
$ws.Range("J102").Value = 2988.889
$ws.Range("K102").Value = 1589
$ws.Range("L102").Value = 2988.889
$ws.Range("M102").Value = 33
$ws.Range("N102").Value = -6232.889

$ws.Range("H132").Value = 2816.8572
$ws.Range("I132").Value = 1702.4
$ws.Range("K132").Value = 5107.200000000001
$ws.Range("M132").Value = -2577.200000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 142859100
$ws.Range("I22").Value = 166667280
$ws.Range("K22").Value = 166667280
$ws.Range("M22").Value = -166666985

$ws.Range("H27").Value = 142859100
$ws.Range("I27").Value = 166667280
$ws.Range("K27").Value = 166667280
$ws.Range("M27").Value = -166667173

$ws.Range("H132").Value = 3205.56
$ws.Range("I132").Value = 2262.5833
$ws.Range("J132").Value = 4076
$ws.Range("K132").Value = 6787.749899999999
$ws.Range("L132").Value = 12228
$ws.Range("M132").Value = -4257.749899999999
$ws.Range("N132").Value = -17288

$ws.Range("H135").Value = 32286
$ws.Range("J135").Value = 32286
$ws.Range("L135").Value = 32286
$ws.Range("N135").Value = -42426

$ws.Range("H136").Value = 1795.25
$ws.Range("I136").Value = 1540.2963
$ws.Range("J136").Value = 2200.1765
$ws.Range("K136").Value = 4620.8889
$ws.Range("L136").Value = 6600.529500000001
$ws.Range("M136").Value = -2070.8889
$ws.Range("N136").Value = -11700.5295

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 13500
$ws.Range("J54").Value = 13500
$ws.Range("L54").Value = 13500
$ws.Range("N54").Value = -14540

$ws.Range("H62").Value = 3728.5715
$ws.Range("I62").Value = 3400
$ws.Range("K62").Value = 3400
$ws.Range("M62").Value = -2776

$ws.Range("H65").Value = 3728.5715
$ws.Range("I65").Value = 3400
$ws.Range("K65").Value = 17000
$ws.Range("M65").Value = -13880

$ws.Range("H126").Value = 2859081.2
$ws.Range("I126").Value = 1558.5714
$ws.Range("K126").Value = 4675.7142
$ws.Range("M126").Value = -2205.7142

$ws.Range("H132").Value = 15632.132
$ws.Range("I132").Value = 1901.52
$ws.Range("J132").Value = 42037.152
$ws.Range("K132").Value = 5704.559999999999
$ws.Range("L132").Value = 126111.456
$ws.Range("M132").Value = -3174.559999999999
$ws.Range("N132").Value = -131171.456

$ws.Range("H136").Value = 2592.158
$ws.Range("I136").Value = 1260.1
$ws.Range("J136").Value = 4072.2222
$ws.Range("K136").Value = 3780.3
$ws.Range("L136").Value = 12216.6666
$ws.Range("M136").Value = -1230.3
$ws.Range("N136").Value = -17316.6666
